# Word COM-interop script.
#
# Summary of the edit (per the commit's xml diff):
#   1. The trailing sentence of the document, "... and see it by
#      yourself!" loses its exclamation mark, which is replaced by a
#      full stop, emitted as its own run: "and see it by yourself" + ".".
#   2. The hidden "_GoBack" bookmark (an empty/collapsed bookmark used
#      by Word to remember the last edit location) is relocated from
#      right after "... depending on their masses," (its position in
#      the original file) to the very end of the document body text,
#      immediately after the new final "." run.
#
# Everything else in the document body is textually unchanged.

$d = $word.ActiveDocument

# --- Step 1: "and see it by yourself!" -> "and see it by yourself." ---
# Locate the exact run of text so we only touch this one occurrence.
$sentence = $d.Range(0, $d.Content.End)
$sentence.Find.Execute("and see it by yourself!", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$exclamationMark = $d.Range($sentence.End - 1, $sentence.End)
$exclamationMark.Text = "."

# Force a run split right after "yourself" / before the new "." so the
# full stop ends up as its own <w:r> (matching the target markup),
# while still inheriting the surrounding run's formatting (en-GB
# language run property) because the split happens in-place on a run
# that already carried it.
$period = $d.Range($sentence.End - 1, $sentence.End)
$period.Font.Bold = $true
$period = $d.Range($sentence.End - 1, $sentence.End)
$period.Font.Bold = $false

# --- Step 2: move the _GoBack bookmark to the end of the document ---
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

$endOfText = $d.Range($sentence.End, $sentence.End)
$d.Bookmarks.Add("_GoBack", $endOfText) | Out-Null
